# Apply the "alpha_zero" nonconvex experiment values (all cases except the 5th)
# to the MitsosBarton2006Ex312 Strong Stationary generator workbook.
#
# All of the affected cells in this workbook hold their numbers as plain TEXT
# (shared-string) values rather than as native Excel numbers, so we must make
# sure the new values keep that same text representation instead of being
# auto-converted into numeric cells by Excel's smart input parsing.

$wb = $excel.ActiveWorkbook

function Set-TextValue($Range, $Text) {
    # Excel will silently reinterpret a numeric-looking string as a real
    # number when assigned directly to .Value. Prefixing with an apostrophe
    # forces Excel to keep it as text (like typing '-1.8 into a cell).
    $Range.Value = "'" + $Text
    # The quote-prefix input leaves a "quote prefix" flag/number format on
    # the cell's style; restore the default "Normal" style so the cell's
    # formatting stays exactly as it was before.
    $Range.Style = "Normal"
}

# NOTE: worksheet names are resolved case-insensitively via Worksheets.Item(name),
# and this workbook has two sheets whose names differ only by case
# ("Vector_bf" vs "Vector_BF"), which would collide if looked up by name.
# Using the (unambiguous) 1-based tab index avoids that problem:
#   1 Funciones_Objetivo
#   2 Restricciones_del_lider
#   3 Restricciones_del_follower
#   4 Punto_modificado
#   5 Vector_bf
#   6 Vector_BF
#   7 Vector_Alpha

function Get-NamedSheet($Index, $ExpectedName) {
    $sheet = $wb.Worksheets.Item($Index)
    if ($sheet.Name -ne $ExpectedName) {
        throw "Expected worksheet #$Index to be '$ExpectedName' but found '$($sheet.Name)'"
    }
    return $sheet
}

# --- Restricciones_del_lider ---------------------------------------------
$wsLider = Get-NamedSheet 2 "Restricciones_del_lider"

Set-TextValue $wsLider.Range("A2") "0.8 - x"
Set-TextValue $wsLider.Range("B2") "-1.8"
Set-TextValue $wsLider.Range("D2") "0.74"

Set-TextValue $wsLider.Range("A3") "-0.8 + x"
Set-TextValue $wsLider.Range("B3") "-0.19999999999999996"
Set-TextValue $wsLider.Range("D3") "0.96"

# --- Restricciones_del_follower -------------------------------------------
$wsFollower = Get-NamedSheet 3 "Restricciones_del_follower"

Set-TextValue $wsFollower.Range("A2") "1.85 - y"
Set-TextValue $wsFollower.Range("B2") "-2.85"
Set-TextValue $wsFollower.Range("D2") "0.76"
Set-TextValue $wsFollower.Range("E2") "8.4"
Set-TextValue $wsFollower.Range("F2") "0.1"

Set-TextValue $wsFollower.Range("A3") "-1.85 + y"
Set-TextValue $wsFollower.Range("B3") "0.8500000000000001"
Set-TextValue $wsFollower.Range("D3") "0.2"
Set-TextValue $wsFollower.Range("E3") "4.5"
Set-TextValue $wsFollower.Range("F3") "2.5"

# --- Punto_modificado -------------------------------------------------------
$wsPunto = Get-NamedSheet 4 "Punto_modificado"

Set-TextValue $wsPunto.Range("A2") "0.8"
Set-TextValue $wsPunto.Range("B2") "1.85"

# --- Vector_bf ---------------------------------------------------------------
$wsBf = Get-NamedSheet 5 "Vector_bf"

Set-TextValue $wsBf.Range("A2") "-9.14325"

# --- Vector_BF -----------------------------------------------------------------
$wsBF = Get-NamedSheet 6 "Vector_BF"

Set-TextValue $wsBF.Range("A2") "-1.07"
Set-TextValue $wsBF.Range("A3") "-33.9"

# Vector_Alpha sheet is intentionally left untouched (the 5th case is excluded).
